$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 393: new positive cases count bumped 92 -> 93 (cumulative column B recalculates itself)
$ws.Range("C393").Value = 93

# Row 407: new positive cases revised 150 -> 149, and one extra-hospital death (M) recorded
$ws.Range("C407").Value = 149
# Column M (and L) are formatted as Text (@); flip to General so the COM
# assignment stores a real number, then restore the Text format so the
# number format itself is unchanged.
$mFmt = $ws.Range("M407").NumberFormat
$ws.Range("M407").NumberFormat = "General"
$ws.Range("M407").Value = 1
$ws.Range("M407").NumberFormat = $mFmt

# Row 409: new positive cases revised 123 -> 127
$ws.Range("C409").Value = 127

# Row 411: new positive cases revised 24 -> 40
$ws.Range("C411").Value = 40

# Row 412: new positive cases revised 11 -> 129
$ws.Range("C412").Value = 129

# Row 413: new day's data filled in (previously a blank placeholder row)
$ws.Range("C413").Value = 9
$ws.Range("E413").Value = 9
$ws.Range("F413").Value = 9
$ws.Range("G413").Value = 38

$lFmt = $ws.Range("L413").NumberFormat
$ws.Range("L413").NumberFormat = "General"
$ws.Range("L413").Value = 0
$ws.Range("L413").NumberFormat = $lFmt

$mFmt413 = $ws.Range("M413").NumberFormat
$ws.Range("M413").NumberFormat = "General"
$ws.Range("M413").Value = 0
$ws.Range("M413").NumberFormat = $mFmt413
